$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42601.976909722223, "Named", 3149, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(42601.988217592596, "Named", 3142, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(42601.98878472222,  "Named", 3066, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(42601.990972222222, "Named", 3243, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(42601.994826388887, "Named", 3171, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(42601.997800925928, "Named", 2850, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40),
    @(42601.99895833333,  "Named", 3073, 120, 5, 2, 0, 100, 0, 3, 2, 60, 40)
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
